# "version with dual gmt tx"
# Insert a new row 69 on Sheet1 (a second GMT tx link entry), pushing the
# former row 69 (DAQ tx link) down to row 70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 69 (the old "DAQ tx" row), which
# shifts it down to row 70 and copies the formatting of the row above
# (row 68, the existing "gmt_tx" row) into the freshly-inserted row 69.
$ws.Rows.Item(69).Insert()

# Match row 68's row height (13.8) on the new row as well.
$ws.Rows.Item(69).RowHeight = 13.8

# Populate the new row 69 with the second "gmt_tx" link entry.
$ws.Range("A69").Value = "tx"
$ws.Range("B69").Value = "1"
$ws.Range("C69").Value = "36"
$ws.Range("D69").Value = "CPPF"
$ws.Range("E69").Value = "Q"
$ws.Range("F69").Value = "125"
$ws.Range("G69").Value = "gmt_tx"
$ws.Range("H69").Value = "1"

# Reflect the author's final view state: scrolled down with C68 selected.
$ws.Range("C68").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
